# Update "想去人数" (interested-attendee count) figures for two events
# on the "展览" and "全部类型" sheets, reflecting newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1718   # 南宁·草莓动漫节: 1716 -> 1718
$ws1.Range("F5").Value = 766    # 南宁·第一届ANE·DACG动漫嘉年华: 765 -> 766

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1718   # 南宁·草莓动漫节: 1716 -> 1718
$ws4.Range("F6").Value = 766    # 南宁·第一届ANE·DACG动漫嘉年华: 765 -> 766
